$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 62
$ws.Range("H62").Value = 67852.94
$ws.Range("I62").Value = 103108.55
$ws.Range("K62").Value = 103108.55
$ws.Range("M62").Value = -102484.55
# Row 65
$ws.Range("H65").Value = 67852.94
$ws.Range("I65").Value = 103108.55
$ws.Range("K65").Value = 515542.75
$ws.Range("M65").Value = -512422.75
# Row 98
$ws.Range("H98").Value = 1257.091
$ws.Range("J98").Value = 567.6667
$ws.Range("L98").Value = 567.6667
$ws.Range("N98").Value = -3563.6667
# Row 122
$ws.Range("H122").Value = 1257.091
$ws.Range("J122").Value = 567.6667
$ws.Range("L122").Value = 1703.0001
$ws.Range("N122").Value = -6603.0001
# Row 137
$ws.Range("H137").Value = 924.63635
$ws.Range("I137").Value = 807.4737
$ws.Range("J137").Value = 1666.6666
$ws.Range("K137").Value = 2422.4211
$ws.Range("L137").Value = 4999.9998
$ws.Range("M137").Value = 127.5789
$ws.Range("N137").Value = -10099.9998

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 1343
$ws.Range("I61").Value = 1011.6
$ws.Range("J61").Value = 3000
$ws.Range("K61").Value = 1011.6
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = -799.6
$ws.Range("N61").Value = -3424
# Row 74
$ws.Range("H74").Value = 1490.0476
$ws.Range("I74").Value = 1627.9286
$ws.Range("J74").Value = 1214.2858
$ws.Range("K74").Value = 1627.9286
$ws.Range("L74").Value = 1214.2858
$ws.Range("M74").Value = -753.9286
$ws.Range("N74").Value = -2962.2858
# Row 77
$ws.Range("H77").Value = 1490.0476
$ws.Range("I77").Value = 1627.9286
$ws.Range("J77").Value = 1214.2858
$ws.Range("K77").Value = 8139.643
$ws.Range("L77").Value = 6071.429
$ws.Range("M77").Value = -3771.643
$ws.Range("N77").Value = -14807.429
# Row 132
$ws.Range("H132").Value = 1951.1923
$ws.Range("I132").Value = 1220.6
$ws.Range("J132").Value = 2407.8125
$ws.Range("K132").Value = 3661.8
$ws.Range("L132").Value = 7223.4375
$ws.Range("M132").Value = -1131.8
$ws.Range("N132").Value = -12283.4375
# Row 136
$ws.Range("H136").Value = 1343
$ws.Range("I136").Value = 1011.6
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 3034.8
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -484.8000000000002
$ws.Range("N136").Value = -14100

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 37912.543
$ws.Range("I134").Value = 3979.7856
$ws.Range("J134").Value = 121746.414
$ws.Range("K134").Value = 11939.3568
$ws.Range("L134").Value = 365239.242
$ws.Range("M134").Value = -9404.356800000001
$ws.Range("N134").Value = -370309.242

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2019.1765
$ws.Range("I31").Value = 2019.1765
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 2019.1765
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -1724.1765
$ws.Range("N31").ClearContents()
# Row 34
$ws.Range("H34").Value = 2019.1765
$ws.Range("I34").Value = 2019.1765
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 2019.1765
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -1817.1765
$ws.Range("N34").ClearContents()
# Row 58
$ws.Range("H58").Value = 3496.639
$ws.Range("I58").Value = 498.1154
$ws.Range("J58").Value = 11292.8
$ws.Range("K58").Value = 498.1154
$ws.Range("L58").Value = 11292.8
$ws.Range("M58").Value = -295.1154
$ws.Range("N58").Value = -11698.8
# Row 86
$ws.Range("H86").Value = 9777.346
$ws.Range("I86").Value = 18770.666
$ws.Range("K86").Value = 18770.666
$ws.Range("M86").Value = -17647.666
# Row 89
$ws.Range("H89").Value = 9777.346
$ws.Range("I89").Value = 18770.666
$ws.Range("K89").Value = 93853.33
$ws.Range("M89").Value = -88237.33
# Row 132
$ws.Range("H132").Value = 1585.4286
$ws.Range("I132").Value = 823
$ws.Range("J132").Value = 2957.8
$ws.Range("K132").Value = 2469
$ws.Range("L132").Value = 8873.400000000001
$ws.Range("M132").Value = 61
$ws.Range("N132").Value = -13933.4
# Row 134
$ws.Range("H134").Value = 2049.7
$ws.Range("I134").Value = 1278.7916
$ws.Range("K134").Value = 3836.3748
$ws.Range("M134").Value = -1301.3748
# Row 136
$ws.Range("H136").Value = 3496.639
$ws.Range("I136").Value = 498.1154
$ws.Range("J136").Value = 11292.8
$ws.Range("K136").Value = 1494.3462
$ws.Range("L136").Value = 33878.39999999999
$ws.Range("M136").Value = 1055.6538
$ws.Range("N136").Value = -38978.39999999999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 776.0714
$ws.Range("J131").Value = 893
$ws.Range("L131").Value = 2679
$ws.Range("N131").Value = -12759

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 62
$ws.Range("H62").Value = 15000
$ws.Range("J62").Value = 15000
$ws.Range("L62").Value = 15000
$ws.Range("N62").Value = -16372
# Row 65
$ws.Range("H65").Value = 15000
$ws.Range("J65").Value = 15000
$ws.Range("L65").Value = 45000
$ws.Range("N65").Value = -51864
# Row 122
$ws.Range("H122").Value = 629961.4399999999
$ws.Range("I122").Value = 1882113.4
$ws.Range("J122").Value = 3885.4285
$ws.Range("K122").Value = 5646340.199999999
$ws.Range("L122").Value = 11656.2855
$ws.Range("M122").Value = -5643890.199999999
$ws.Range("N122").Value = -16556.2855
# Row 126
$ws.Range("H126").Value = 3046.8696
$ws.Range("I126").Value = 3126.8235
$ws.Range("J126").Value = 2820.3333
$ws.Range("K126").Value = 9380.470499999999
$ws.Range("L126").Value = 8460.999899999999
$ws.Range("M126").Value = -6910.470499999999
$ws.Range("N126").Value = -13400.9999
# Row 132
$ws.Range("H132").Value = 3859.647
$ws.Range("I132").Value = 4101.778
$ws.Range("J132").Value = 3587.25
$ws.Range("K132").Value = 12305.334
$ws.Range("L132").Value = 10761.75
$ws.Range("M132").Value = -9775.334000000001
$ws.Range("N132").Value = -15821.75

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 122
$ws.Range("H122").Value = 3587.2546
$ws.Range("I122").Value = 4606.3794
$ws.Range("K122").Value = 13819.1382
$ws.Range("M122").Value = -11369.1382
# Row 132
$ws.Range("H132").Value = 3062.75
$ws.Range("I132").Value = 2743.4517
$ws.Range("K132").Value = 8230.355100000001
$ws.Range("M132").Value = -5700.355100000001
# Row 136
$ws.Range("H136").Value = 5149.087
$ws.Range("I136").Value = 1499.6
$ws.Range("J136").Value = 11991.875
$ws.Range("K136").Value = 4498.799999999999
$ws.Range("L136").Value = 35975.625
$ws.Range("M136").Value = -1948.799999999999
$ws.Range("N136").Value = -41075.625

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 1189.1428
$ws.Range("I122").Value = 1113.4546
$ws.Range("K122").Value = 3340.3638
$ws.Range("M122").Value = -890.3638000000001
# Row 132
$ws.Range("H132").Value = 1647.3334
$ws.Range("I132").Value = 1398.68
$ws.Range("J132").Value = 2890.6
$ws.Range("K132").Value = 4196.04
$ws.Range("L132").Value = 8671.799999999999
$ws.Range("M132").Value = -1666.04
$ws.Range("N132").Value = -13731.8
# Row 136
$ws.Range("H136").Value = 1371.6586
$ws.Range("I136").Value = 1413.0294
$ws.Range("J136").Value = 1170.7142
$ws.Range("K136").Value = 4239.0882
$ws.Range("L136").Value = 1170.7142
$ws.Range("M136").Value = -1689.0882
$ws.Range("N136").Value = -8612.142599999999
